# Applies the "Updated cryptos list" data refresh to Price (D) and Volume(1h) (E) columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.149.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.480.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.33%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.479.98"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("E10").Value = "  +2.97%  "
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.929.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.970.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.486.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.88%  "
$ws.Range("E20").Value = "  -3.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "350.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.51%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.607.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.35%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0911"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "505.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.08%  "
$ws.Range("E34").Value = "  -0.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.27%  "
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.61%  "
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.329"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.10%  "
$ws.Range("E43").Value = "  +1.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "143.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₆0265"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.515"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0738"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("E51").Value = "  -0.93%  "

Write-Output "Applied cryptos update"
